$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 17-21 were type "专职单" / rows 22-26 were type "暑假单".
# Both groups are being reclassified as "好单" (a "good order"), and since
# those two strings become unused, they disappear from the shared-string
# table on save.
for ($r = 17; $r -le 26; $r++) {
    $ws.Cells.Item($r, 3).Value = "好单"
}

# Rows 18-26 also get their "creation date" (column E) bumped from
# 2025/6/13 to 2025/6/14; row 17 keeps its original date.
for ($r = 18; $r -le 26; $r++) {
    $ws.Cells.Item($r, 5).Value = "2025/6/14"
}

# The author's Excel window was resized (e.g. maximized on a bigger
# screen) between edits.
$excel.ActiveWindow.Width = 22188
$excel.ActiveWindow.Height = 9767

# The view had scrolled so D3 was the top-left visible cell with D5
# selected; afterwards the sheet is scrolled back to the top (no
# topLeftCell override) and the selection sits on G25.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G25").Select()
